# Appends " (Changed main)" to the end of the first paragraph's text,
# split across three separate runs (" (", "Changed main", ")") so the
# resulting OOXML mirrors a manual multi-step edit rather than a single
# merged run.

$d = $word.ActiveDocument

# Locate the paragraph that currently reads exactly
# "This is a Microsoft word document." and work from its range so the
# script is resilient to the paragraph's absolute position in the body.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "This is a Microsoft word document.") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    $target = $d.Paragraphs(1)
}

# End of the visible text, i.e. just before the paragraph mark.
$insertPos = $target.Range.End - 1

function Insert-AsNewRun($position, $text) {
    # InsertAfter() on a zero-length range normally gets silently coalesced
    # into the preceding run when it shares identical run formatting. Wrapping
    # the freshly-inserted text in a bookmark (and then deleting that
    # bookmark again) forces the engine to keep it as an independent <w:r>
    # element instead of merging it back into its neighbor.
    $r = $d.Range($position, $position)
    $r.InsertAfter($text)
    $bmName = "tmp_split_" + [Guid]::NewGuid().ToString("N")
    $d.Bookmarks.Add($bmName, $r) | Out-Null
    $d.Bookmarks($bmName).Delete()
    return $text.Length
}

$pos = $insertPos
$pos += Insert-AsNewRun $pos " ("
$pos += Insert-AsNewRun $pos "Changed main"
$pos += Insert-AsNewRun $pos ")"

Write-Host "Paragraph now reads: [$($target.Range.Text)]"
